# "Generate Report for Handoff" - refresh localization-status report with
# a new pair of handoff files and updated status/timestamps.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$file1 = "5a04021b-a533-40d8-9da6-7aaf68baea15.md"
$file2 = "fffffcdd500d-9591-4e41-ae7a-bcfc200ac221.md"
$status = "Ready for handoff"
$genDate = "2016-08-25 00:59:06"
$zhHandoffFile = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.zh-cn.xlf"
$deHandoffFile = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.de-de.xlf"
$zhHandoffDate = "2016-08-25 00:58:57"
$nullDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws1.Range("A2").Value = $file1
$ws1.Range("A3").Value = $file2
$ws1.Range("E2").Value = $status
$ws1.Range("F2").Value = $status
$ws1.Range("E3").Value = $status
$ws1.Range("F3").Value = $status
$ws1.Range("G2").Value = $genDate
$ws1.Range("G3").Value = $genDate

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9499c0bded293588ae588fff3391a4ab6671c0b/e2e/1c0696b1-d9a5-4d71-b2b8-029e60c0a26b.md", "", "", "e2e\" + $file1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9499c0bded293588ae588fff3391a4ab6671c0b/e2e/3019a9a2-7081-43f1-ba05-ee32a65e7bf0.md", "", "", "e2e\" + $file2) | Out-Null

$ws1.Columns("E:F").ColumnWidth = 16.25

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2.Range("A2").Value = $file1
$ws2.Range("A3").Value = $file2
$ws2.Range("C2").Value = $status
$ws2.Range("C3").Value = $status
$ws2.Range("F3").Value = "True"
$ws2.Range("G2").Value = $zhHandoffFile
$ws2.Range("G3").Value = $zhHandoffFile
$ws2.Range("H2").Value = $zhHandoffDate
$ws2.Range("H3").Value = $zhHandoffDate
$ws2.Range("I2").Value = ""
$ws2.Range("I2").ClearFormats()
$ws2.Range("I3").Value = ""
$ws2.Range("I3").ClearFormats()
$ws2.Range("J2").Value = ""
$ws2.Range("J3").Value = ""
$ws2.Range("K2").Value = $nullDate
$ws2.Range("K3").Value = $nullDate

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9499c0bded293588ae588fff3391a4ab6671c0b/e2e/1c0696b1-d9a5-4d71-b2b8-029e60c0a26b.md", "", "", $file1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9499c0bded293588ae588fff3391a4ab6671c0b/e2e/3019a9a2-7081-43f1-ba05-ee32a65e7bf0.md", "", "", $file2) | Out-Null

$ws2.Columns("C").ColumnWidth = 16.25
$ws2.Columns("I").ColumnWidth = 17.77
$ws2.Columns("J").ColumnWidth = 20.77

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3.Range("A2").Value = $file1
$ws3.Range("A3").Value = $file2
$ws3.Range("C2").Value = $status
$ws3.Range("C3").Value = $status
$ws3.Range("F3").Value = "True"
$ws3.Range("G2").Value = $deHandoffFile
$ws3.Range("G3").Value = $deHandoffFile
$ws3.Range("H2").Value = $genDate
$ws3.Range("H3").Value = $genDate
$ws3.Range("I2").Value = ""
$ws3.Range("I2").ClearFormats()
$ws3.Range("I3").Value = ""
$ws3.Range("I3").ClearFormats()
$ws3.Range("J2").Value = ""
$ws3.Range("J3").Value = ""
$ws3.Range("K2").Value = $nullDate
$ws3.Range("K3").Value = $nullDate

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9499c0bded293588ae588fff3391a4ab6671c0b/e2e/1c0696b1-d9a5-4d71-b2b8-029e60c0a26b.md", "", "", $file1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9499c0bded293588ae588fff3391a4ab6671c0b/e2e/3019a9a2-7081-43f1-ba05-ee32a65e7bf0.md", "", "", $file2) | Out-Null

$ws3.Columns("C").ColumnWidth = 16.25
$ws3.Columns("I").ColumnWidth = 17.77
$ws3.Columns("J").ColumnWidth = 20.77
